$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - updated "want to go" counts (column F)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3084
$ws.Range("F5").Value = 1686
$ws.Range("F6").Value = 2093
$ws.Range("F9").Value = 894
$ws.Range("F11").Value = 219
$ws.Range("F12").Value = 441
$ws.Range("F13").Value = 1146
$ws.Range("F17").Value = 7491
$ws.Range("F18").Value = 315
$ws.Range("F25").Value = 74
$ws.Range("F26").Value = 1127
$ws.Range("F27").Value = 969
$ws.Range("F29").Value = 1550
$ws.Range("F31").Value = 1139
$ws.Range("F33").Value = 478
$ws.Range("F36").Value = 263
$ws.Range("F39").Value = 311
$ws.Range("F41").Value = 206

# Sheet "演出" (shows)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 21

# Sheet "全部类型" (all types, combined listing)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 21
$ws.Range("F7").Value = 3084
$ws.Range("F8").Value = 1686
$ws.Range("F9").Value = 2093
$ws.Range("F12").Value = 894
$ws.Range("F15").Value = 219
$ws.Range("F16").Value = 441
$ws.Range("F17").Value = 1146
$ws.Range("F21").Value = 7491
$ws.Range("F22").Value = 315
$ws.Range("F30").Value = 74
$ws.Range("F31").Value = 1127
$ws.Range("F32").Value = 969
$ws.Range("F34").Value = 1550
$ws.Range("F36").Value = 1139
$ws.Range("F38").Value = 478
$ws.Range("F41").Value = 263
$ws.Range("F44").Value = 311
$ws.Range("F49").Value = 206
